$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the values among rows 2, 3 and 6 (row2 <- row3 <- row6 <- row2)
$oldA2 = $ws.Range("A2").Value2
$oldB2 = $ws.Range("B2").Value2
$oldA3 = $ws.Range("A3").Value2
$oldB3 = $ws.Range("B3").Value2
$oldA6 = $ws.Range("A6").Value2
$oldB6 = $ws.Range("B6").Value2

$ws.Range("A2").Value = $oldA3
$ws.Range("B2").Value = $oldB3

$ws.Range("A3").Value = $oldA6
$ws.Range("B3").Value = $oldB6

$ws.Range("A6").Value = $oldA2
$ws.Range("B6").Value = $oldB2

# Update the selected cell in the sheet view
$ws.Range("K22").Select()
